$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.215.55"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.647.56"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "1.877.83"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "1.667.42"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "27.204.25"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "0.0₃0743"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0507"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "1.267.04"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.845"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("E42").Value = "  +4.62%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "1.788.59"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("E51").Value = "  -0.60%  "
